$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hyperparameter Optimization")

# Row 17
$ws.Range("F17").Value = 0.038574970441669201
$ws.Range("G17").Value = "Adam"
$ws.Range("H17").Value = "relu"
$ws.Range("I17").Value = "256"
$ws.Range("J17").Value = 0.193868361915666
$ws.Range("K17").Value = "50"
$ws.Range("L17").Value = 0.064655070006847296
$ws.Range("M17").Value = 0.0131740950694645
$ws.Range("P17").Value = "[0.051007986068725586, 0.04874471575021744, 0.07474925369024277, 0.08266618847846985, 0.06610720604658127]"

# Row 18
$ws.Range("F18").Value = 0.0036702828800194401
$ws.Range("G18").Value = "Adam"
$ws.Range("H18").Value = "tanh"
$ws.Range("I18").Value = "256"
$ws.Range("J18").Value = 0.17556384197400199
$ws.Range("K18").Value = "31"
$ws.Range("L18").Value = 0.064543188363313603
$ws.Range("M18").Value = 0.0114796169486963
$ws.Range("P18").Value = "[0.0517902709543705, 0.05279034376144409, 0.07151271402835846, 0.08213762193918228, 0.06448499113321304]"

# Row 19
$ws.Range("F19").Value = 0.0081780459091270706
$ws.Range("G19").Value = "SGD"
$ws.Range("H19").Value = "sigmoid"
$ws.Range("I19").Value = "512"
$ws.Range("J19").Value = 0.14852115505240601
$ws.Range("K19").Value = "16"
$ws.Range("L19").Value = 0.063691043853759705
$ws.Range("M19").Value = 0.0118308326479514
$ws.Range("P19").Value = "[0.052902087569236755, 0.049781009554862976, 0.07142113149166107, 0.08182501047849655, 0.06252598017454147]"

# Row 57
$ws.Range("E57").Value = "1024"
$ws.Range("F57").Value = 0.060770659627901798
$ws.Range("G57").Value = "SGD"
$ws.Range("H57").Value = "tanh"
$ws.Range("I57").Value = "256"
$ws.Range("J57").Value = 0.28693073271555403
$ws.Range("K57").Value = "76"
$ws.Range("L57").Value = 0.051908043771982097
$ws.Range("M57").Value = 0.0088052118526808394
$ws.Range("P57").Value = "[0.044655557721853256, 0.043497055768966675, 0.056257907301187515, 0.06708838045597076, 0.048041317611932755]"

# Row 58
$ws.Range("E58").Value = "1024"
$ws.Range("F58").Value = 0.21903465162979599
$ws.Range("G58").Value = "SGD"
$ws.Range("H58").Value = "relu"
$ws.Range("I58").Value = "1024"
$ws.Range("J58").Value = 0.49088846084547999
$ws.Range("K58").Value = "87"
$ws.Range("L58").Value = 0.058136474341154099
$ws.Range("M58").Value = 0.011278367874087401
$ws.Range("P58").Value = "[0.04046407714486122, 0.05292431265115738, 0.06083891540765762, 0.07476838678121567, 0.0616866797208786]"

# Row 59
$ws.Range("E59").Value = "1024"
$ws.Range("F59").Value = 0.097276147553062797
$ws.Range("G59").Value = "SGD"
$ws.Range("H59").Value = "relu"
$ws.Range("I59").Value = "1024"
$ws.Range("J59").Value = 0.11844994511316401
$ws.Range("K59").Value = "92"
$ws.Range("L59").Value = 0.0569791264832019
$ws.Range("M59").Value = 0.0092513184450840492
$ws.Range("P59").Value = "[0.04326264560222626, 0.053425583988428116, 0.057177986949682236, 0.07189487665891647, 0.05913453921675682]"

# Update the active sheet's view/selection to match the saved workbook state
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("F20").Select()

$excel.Calculate()
